$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for LavieBio / P22, right after the AgPlenus block (old row 7),
# pushing the rest of the LavieBio rows (and everything below) down by one.
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "LavieBio"

# Rename the team/company associated with a handful of existing budgets
# (the budget codes themselves stay put, only column A changes).
$ws.Range("A29").Value = "Chempass"
$ws.Range("A30").Value = "MicroBoost"
$ws.Range("A31").Value = "MicroBoost"
$ws.Range("A32").Value = "Upkeep"

# Set the new budget code last so the shared-string table grows in the same
# order as the source edit (Chempass, MicroBoost, P22).
$ws.Range("B7").Value = "P22"

# Restore the view: zoom to 85% and move the active selection to G13.
$excel.ActiveWindow.Zoom = 85
$ws.Range("G13").Select()
